$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: replace the repeated "*" symbol with a distinct symbol per row
# (row 2 keeps "*", rows 3-6 get new symbols).
$ws.Range("C3").Value = "#"
$ws.Range("C4").Value = "+"
$ws.Range("C5").Value = "%"
$ws.Range("C6").Value = "^"

# Column A: row 6 value changes from 2354 to 235
$ws.Range("A6").Value = 235

# Move the active selection from D6 to C6
$ws.Range("C6").Select()
